$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force a literal text value (preserving leading zeros, etc.) while
    # keeping the cell's original (default) style - mirrors typing an
    # apostrophe-prefixed entry in Excel, then clearing the quote-prefix
    # style marker back to Normal so no stray formatting sticks around.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("A2") "000001"
$ws.Range("B2").Value = "AAA01"
Set-TextValue $ws.Range("C2") "11/11/2024 10:20:00"
Set-TextValue $ws.Range("D2") "11/11/2024 10:21:00"
$ws.Range("E2").Value = "Moto"
$ws.Range("F2").Value = 1500
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1500

# Row 3
Set-TextValue $ws.Range("A3") "000002"
$ws.Range("B3").Value = "BBB01"
Set-TextValue $ws.Range("C3") "11/11/2024 10:20:00"
Set-TextValue $ws.Range("D3") "11/11/2024 10:25:00"
$ws.Range("E3").Value = "Moto"
$ws.Range("F3").Value = 1500
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1500

# Row 4
Set-TextValue $ws.Range("A4") "000003"
$ws.Range("B4").Value = "AAA01"
Set-TextValue $ws.Range("C4") "11/11/2024 10:23:00"
Set-TextValue $ws.Range("D4") "11/11/2024 10:25:00"
$ws.Range("E4").Value = "Moto"
$ws.Range("F4").Value = 1500
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1500

# Row 5
Set-TextValue $ws.Range("A5") "000004"
$ws.Range("B5").Value = "BBB01"
Set-TextValue $ws.Range("C5") "11/11/2024 10:27:00"
Set-TextValue $ws.Range("D5") "11/11/2024 10:28:00"
$ws.Range("E5").Value = "Moto"
$ws.Range("F5").Value = 1500
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1500

# Row 6
Set-TextValue $ws.Range("A6") "000005"
$ws.Range("B6").Value = "CCC01"
Set-TextValue $ws.Range("C6") "11/11/2024 10:30:00"
Set-TextValue $ws.Range("D6") "11/11/2024 10:31:00"
$ws.Range("E6").Value = "Moto"
$ws.Range("F6").Value = 1500
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1500

# Row 7
Set-TextValue $ws.Range("A7") "000006"
$ws.Range("B7").Value = "DDD01"
Set-TextValue $ws.Range("C7") "11/11/2024 10:44:00"
Set-TextValue $ws.Range("D7") "11/11/2024 10:45:00"
$ws.Range("E7").Value = "Moto"
$ws.Range("F7").Value = 1500
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1500

# Row 8
Set-TextValue $ws.Range("A8") "000007"
$ws.Range("B8").Value = "AAA01"
Set-TextValue $ws.Range("C8") "11/11/2024 11:09:00"
Set-TextValue $ws.Range("D8") "11/11/2024 11:10:00"
$ws.Range("E8").Value = "Moto"
$ws.Range("F8").Value = 1500
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1500

# Row 9
Set-TextValue $ws.Range("A9") "000008"
$ws.Range("B9").Value = "BBB01"
Set-TextValue $ws.Range("C9") "10/11/2024 10:00:00"
Set-TextValue $ws.Range("D9") "11/11/2024 12:36:00"
$ws.Range("E9").Value = "Moto"
$ws.Range("F9").Value = 7000
Set-TextValue $ws.Range("G9") "26:36"
$ws.Range("H9").Value = 18500

# Row 10 is removed entirely from the register (was 000084 / III01)
$ws.Range("A10:H10").Delete()
